$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a new paragraph right after the target paragraph
    $target.Range.InsertParagraphAfter()

    # The newly created paragraph is the one following $target
    $newPara = $target.Next()

    # Set its style to List Bullet and fill in the text
    $newPara.Style = $d.Styles.Item("List Bullet")
    $newPara.Range.Text = "7043088 - Ana Karine Furtado de Carvalho"
}
